# Update NATMI LR-pair sheet with newly recalculated TPM values.
# Sending cluster changes from MuSCs to ECs; only 2 data rows remain
# (previously 4), and the numeric columns are recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra data rows (old rows 4 and 5) first, from the bottom up
# so row indices for the remaining rows are not disturbed.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Row 2: ECs -> Efna3 -> Epha7 -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna3"
$ws.Range("C2").Value = "Epha7"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02551366666666667
$ws.Range("H2").Value = 0.076541
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2421403333333333
$ws.Range("N2").Value = 0.726421
$ws.Range("O2").Value = 0.7386057795451564
$ws.Range("P2").Value = 0.8091041635804498
$ws.Range("Q2").Value = 0.006177887751222222
$ws.Range("R2").Value = 0.055600989761
$ws.Range("S2").Value = 0.7386057795451564
$ws.Range("T2").Value = 0.8091041635804498

# Row 3: ECs -> Efna3 -> Epha7 -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna3"
$ws.Range("C3").Value = "Epha7"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02551366666666667
$ws.Range("H3").Value = 0.076541
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.08569400000000001
$ws.Range("N3").Value = 0.171388
$ws.Range("O3").Value = 0.2613942204548436
$ws.Range("P3").Value = 0.1908958364195503
$ws.Range("Q3").Value = 0.002186368151333334
$ws.Range("R3").Value = 0.013118208908
$ws.Range("S3").Value = 0.2613942204548436
$ws.Range("T3").Value = 0.1908958364195503
